# Add two new columns, "I0" (col I) and "IF" (col J), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - match the style already used by the other header
# cells (B1:H1): bold font, thin border, centered horizontal/top vertical
# alignment.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data rows 2-74: I0 and IF values, one pair per data row.
$i0Values = @(3,7,7,6,8,7,4,6,6,8,8,7,7,7,6,7,9,7,6,8,7,9,6,9,7,6,6,10,6,7,8,8,8,8,6,8,9,7,7,10,7,7,8,8,8,8,8,7,8,8,9,7,7,9,8,8,7,8,6,8,7,9,8,8,7,6,5,7,7,6,5,8,6)
$ifValues = @(3,7,7,7,8,7,5,6,6,8,8,7,7,7,6,7,9,7,6,8,7,9,7,9,7,6,6,10,6,7,8,9,9,8,7,8,9,7,8,10,7,8,8,8,8,8,8,7,8,8,9,7,8,9,8,8,7,8,6,8,7,9,8,8,7,6,5,7,7,6,5,8,6)

$startRow = 2
for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
